# Add a new JST connector part to the parts list (row 17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D17 first so the new shared-string for the Amazon URL is registered
# before the part-name string, matching the expected shared-strings order.
$ws.Range("D17").Value = "https://www.amazon.com/eBoot-Pairs-Electrical-Female-Connector/dp/B06WGM9W7S/ref=sr_1_1?ie=UTF8&qid=1531752427&sr=8-1&keywords=20+pair+micro+jst+1.25"
$ws.Range("A17").Value = "jst 2 pin 1.25mm electrical female plug"
